$d = $word.ActiveDocument
$brk = [char]11

# 1. Delete the "LOB1012 -  Estatística  (Requisito)" run, including its line break.
$rngLob = $d.Content
$foundLob = $rngLob.Find.Execute("LOB1012 -  Estatística  (Requisito)")
if ($foundLob) {
    $delLobRange = $d.Range($rngLob.Start, $rngLob.End + 1)
    $delLobRange.Delete()
}

# 2. Replace the Química requisitos block (3 consecutive lines) with the new 3 lines.
#    Delete the whole old block first, then insert each new line as its own run
#    (mirrors the diff, which keeps each course on its own <w:r>).
$rngQuim = $d.Content
$line1 = "LOQ4031 -  Química Geral I  (Requisito)"
$line2 = "LOQ4073 -  Química Geral II  (Requisito)"
$line3 = "LOQ4095 -  Química Geral Experimental  (Requisito)"
$pattern = $line1 + $brk + $line2 + $brk + $line3 + $brk
$foundQuim = $rngQuim.Find.Execute($pattern, $true, $false, $false)
if ($foundQuim) {
    $insertPoint = $rngQuim.Start
    $delQuimRange = $d.Range($rngQuim.Start, $rngQuim.End)
    $delQuimRange.Delete()

    $p1 = $d.Range($insertPoint, $insertPoint)
    $p1.InsertAfter("LOQ4095 -  Química Geral Experimental  (Requisito)" + $brk)

    $p2 = $d.Range($p1.End, $p1.End)
    $p2.InsertAfter("LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + $brk)

    $p3 = $d.Range($p2.End, $p2.End)
    $p3.InsertAfter("LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)" + $brk)
}
